# ---------------------------------------------------------------------------
# paises.xlsx - "Datos actualizados" data refresh (5 Sep 2020, 14:13 -> 15:30)
#
# The source feed re-ranked a handful of countries (Suecia, Paises Bajos,
# Libia, Gambia) and refreshed the case / recovered / critical / death
# counters for many rows. Because the sheet is ordered by total cases, the
# re-rank shows up as new country names for several rows in column A, while
# the numeric columns B:H simply carry the latest reported figures. The
# banner timestamp in A1 is also bumped forward.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados a 5 de Septiembre de 2020 a las ..." banner
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 15:30"

# Row 4
$ws.Range("B4").Value = 6390240
$ws.Range("C4").Value = 1183
$ws.Range("D4").Value = 3636284
$ws.Range("E4").Value = 2561810
# Row 13
$ws.Range("D13").Value = 340381
$ws.Range("E13").Value = 111816
$ws.Range("G13").Value = 62
$ws.Range("H13").Value = 9685
# Row 18
$ws.Range("B18").Value = 319932
$ws.Range("C18").Value = 791
$ws.Range("D18").Value = 295842
$ws.Range("E18").Value = 20041
$ws.Range("G18").Value = 34
$ws.Range("H18").Value = 4049
# Row 24
$ws.Range("B24").Value = 250554
$ws.Range("C24").Value = 273
$ws.Range("E24").Value = 14945
# Row 42  # -> Suecia
$ws.Range("A42").Value = "Suecia"
$ws.Range("B42").Value = 84985
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("H42").Value = 5835
# Row 43  # -> Guatemala
$ws.Range("A43").Value = "Guatemala"
$ws.Range("B43").Value = 77040
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 65029
$ws.Range("E43").Value = 9186
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 2825
# Row 44
$ws.Range("B44").Value = 73862
$ws.Range("C44").Value = 654
$ws.Range("G44").Value = 4
$ws.Range("H44").Value = 6241
# Row 45  # -> Emiratos Arabes Unidos
$ws.Range("A45").Value = "Emiratos Arabes Unidos"
$ws.Range("B45").Value = 73471
$ws.Range("C45").Value = 705
$ws.Range("D45").Value = 63652
$ws.Range("E45").Value = 9431
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 388
# Row 46  # -> Bielorrusia
$ws.Range("A46").Value = "Bielorrusia"
$ws.Range("B46").Value = 72663
$ws.Range("C46").Value = 178
$ws.Range("D46").Value = 71843
$ws.Range("E46").Value = 115
$ws.Range("G46").Value = 4
$ws.Range("H46").Value = 705
# Row 47  # -> Polonia
$ws.Range("A47").Value = "Polonia"
$ws.Range("B47").Value = 70387
$ws.Range("C47").Value = 567
$ws.Range("D47").Value = 52346
$ws.Range("E47").Value = 15928
$ws.Range("G47").Value = 13
$ws.Range("H47").Value = 2113
# Row 48  # -> Japon
$ws.Range("A48").Value = "Japon"
$ws.Range("B48").Value = 70268
$ws.Range("D48").Value = 60417
$ws.Range("E48").Value = 8521
$ws.Range("H48").Value = 1330
# Row 49  # -> Marruecos
$ws.Range("A49").Value = "Marruecos"
$ws.Range("B49").Value = 68605
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 52483
$ws.Range("E49").Value = 14830
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 1292
# Row 50  # -> Honduras
$ws.Range("A50").Value = "Honduras"
$ws.Range("B50").Value = 63798
$ws.Range("C50").Value = 640
$ws.Range("D50").Value = 12347
$ws.Range("E50").Value = 49467
$ws.Range("G50").Value = 30
$ws.Range("H50").Value = 1984
# Row 51  # -> Portugal
$ws.Range("A51").Value = "Portugal"
$ws.Range("B51").Value = 59457
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 42576
$ws.Range("E51").Value = 15048
$ws.Range("H51").Value = 1833
# Row 52  # -> Singapur
$ws.Range("A52").Value = "Singapur"
$ws.Range("B52").Value = 56982
$ws.Range("C52").Value = 34
$ws.Range("D52").Value = 56174
$ws.Range("E52").Value = 781
$ws.Range("H52").Value = 27
# Row 53  # -> Etiopia
$ws.Range("A53").Value = "Etiopia"
$ws.Range("B53").Value = 56516
$ws.Range("D53").Value = 20612
$ws.Range("E53").Value = 35024
$ws.Range("H53").Value = 880
# Row 54  # -> Nigeria
$ws.Range("A54").Value = "Nigeria"
$ws.Range("B54").Value = 54743
$ws.Range("D54").Value = 42816
$ws.Range("E54").Value = 10876
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 1051
# Row 55  # -> Barein
$ws.Range("A55").Value = "Barein"
$ws.Range("B55").Value = 54095
$ws.Range("D55").Value = 50323
$ws.Range("E55").Value = 3576
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 196
# Row 56  # -> Venezuela
$ws.Range("A56").Value = "Venezuela"
$ws.Range("B56").Value = 50973
$ws.Range("D56").Value = 41249
$ws.Range("E56").Value = 9312
$ws.Range("H56").Value = 412
# Row 57  # -> Argelia
$ws.Range("A57").Value = "Argelia"
$ws.Range("B57").Value = 45773
$ws.Range("D57").Value = 32259
$ws.Range("E57").Value = 11975
$ws.Range("H57").Value = 1539
# Row 58  # -> Costa Rica
$ws.Range("A58").Value = "Costa Rica"
$ws.Range("B58").Value = 45680
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 18053
$ws.Range("E58").Value = 27158
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 469
# Row 59  # -> Nepal
$ws.Range("A59").Value = "Nepal"
$ws.Range("B59").Value = 45277
$ws.Range("C59").Value = 1041
$ws.Range("D59").Value = 27127
$ws.Range("E59").Value = 17870
$ws.Range("G59").Value = 9
$ws.Range("H59").Value = 280
# Row 60  # -> Ghana
$ws.Range("A60").Value = "Ghana"
$ws.Range("B60").Value = 44777
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 43693
$ws.Range("E60").Value = 801
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 283
# Row 61  # -> Armenia
$ws.Range("A61").Value = "Armenia"
$ws.Range("B61").Value = 44649
$ws.Range("C61").Value = 188
$ws.Range("D61").Value = 39823
$ws.Range("E61").Value = 3931
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 895
# Row 62  # -> Kirguistan
$ws.Range("A62").Value = "Kirguistan"
$ws.Range("B62").Value = 44293
$ws.Range("C62").Value = 94
$ws.Range("D62").Value = 39599
$ws.Range("E62").Value = 3634
$ws.Range("H62").Value = 1060
# Row 63  # -> Suiza
$ws.Range("A63").Value = "Suiza"
$ws.Range("B63").Value = 43957
$ws.Range("C63").Value = 425
$ws.Range("D63").Value = 36500
$ws.Range("E63").Value = 5444
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 2013
# Row 64  # -> Uzbekistan
$ws.Range("A64").Value = "Uzbekistan"
$ws.Range("B64").Value = 43075
$ws.Range("C64").Value = 77
$ws.Range("D64").Value = 40432
$ws.Range("E64").Value = 2304
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 339
# Row 65  # -> Moldavia
$ws.Range("A65").Value = "Moldavia"
$ws.Range("B65").Value = 38906
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 27017
$ws.Range("E65").Value = 10842
$ws.Range("H65").Value = 1047
# Row 66  # -> Afganistan
$ws.Range("A66").Value = "Afganistan"
$ws.Range("B66").Value = 38324
$ws.Range("C66").Value = 20
$ws.Range("D66").Value = 30082
$ws.Range("E66").Value = 6833
$ws.Range("H66").Value = 1409
# Row 67  # -> Azerbaiyan
$ws.Range("A67").Value = "Azerbaiyan"
$ws.Range("B67").Value = 37192
$ws.Range("C67").Value = 161
$ws.Range("D67").Value = 34565
$ws.Range("E67").Value = 2082
$ws.Range("G67").Value = 2
$ws.Range("H67").Value = 545
# Row 68  # -> Kenia
$ws.Range("A68").Value = "Kenia"
$ws.Range("B68").Value = 34884
$ws.Range("D68").Value = 21059
$ws.Range("E68").Value = 13236
$ws.Range("H68").Value = 589
# Row 69  # -> Serbia
$ws.Range("A69").Value = "Serbia"
$ws.Range("B69").Value = 31849
$ws.Range("C69").Value = 77
$ws.Range("D69").Value = 30529
$ws.Range("E69").Value = 597
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 723
# Row 70  # -> Irlanda
$ws.Range("A70").Value = "Irlanda"
$ws.Range("B70").Value = 29303
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 23364
$ws.Range("E70").Value = 4162
$ws.Range("H70").Value = 1777
# Row 71  # -> Austria
$ws.Range("A71").Value = "Austria"
$ws.Range("B71").Value = 29087
$ws.Range("C71").Value = 358
$ws.Range("D71").Value = 24828
$ws.Range("E71").Value = 3524
$ws.Range("H71").Value = 735
# Row 72  # -> Chequia
$ws.Range("A72").Value = "Chequia"
$ws.Range("B72").Value = 27249
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 19027
$ws.Range("E72").Value = 7793
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 429
# Row 73  # -> Australia
$ws.Range("A73").Value = "Australia"
$ws.Range("C73").Value = 71
$ws.Range("D73").Value = 22330
$ws.Range("E73").Value = 3129
$ws.Range("G73").Value = 11
$ws.Range("H73").Value = 748
# Row 82
$ws.Range("B82").Value = 17736
$ws.Range("C82").Value = 189
$ws.Range("D82").Value = 15671
$ws.Range("E82").Value = 1438
# Row 83  # -> Libia
$ws.Range("A83").Value = "Libia"
$ws.Range("B83").Value = 17094
$ws.Range("C83").Value = 649
$ws.Range("D83").Value = 2025
$ws.Range("E83").Value = 14797
$ws.Range("G83").Value = 10
$ws.Range("H83").Value = 272
# Row 84  # -> Bulgaria
$ws.Range("A84").Value = "Bulgaria"
$ws.Range("B84").Value = 16954
$ws.Range("D84").Value = 12046
$ws.Range("E84").Value = 4243
$ws.Range("H84").Value = 665
# Row 98
$ws.Range("B98").Value = 8757
$ws.Range("C98").Value = 33
$ws.Range("D98").Value = 7549
$ws.Range("E98").Value = 1138
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 70
# Row 128  # -> Gambia
$ws.Range("A128").Value = "Gambia"
$ws.Range("B128").Value = 3150
$ws.Range("C128").Value = 30
$ws.Range("D128").Value = 1315
$ws.Range("E128").Value = 1736
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 99
# Row 129  # -> Eslovenia
$ws.Range("A129").Value = "Eslovenia"
$ws.Range("B129").Value = 3122
$ws.Range("C129").Value = 43
$ws.Range("D129").Value = 2483
$ws.Range("E129").Value = 504
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 135
# Row 142
$ws.Range("D142").Value = 1793
$ws.Range("E142").Value = 361
